# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (col E) / date (col F) rows for the second worker
# (ERICK ENRIQUE VALENZUELA ESCOBAR, rows 17-36) were re-keyed: the periods
# were entered in descending order (2311 down to 2204); this update reverses
# them into ascending order (2204 up to 2311), carrying each row's matching
# date along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 17
$lastRow = 36
$count = $lastRow - $firstRow + 1

$periods = @()
$dates = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value2
    $dates += $ws.Cells.Item($r, 6).Value2
}

# Reverse via index slice -- [array]::Reverse() does not mutate in place here.
$periodsRev = $periods[($count - 1)..0]
$datesRev = $dates[($count - 1)..0]

for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 5).Value = $periodsRev[$i]
    $ws.Cells.Item($r, 6).Value = $datesRev[$i]
}
